# AzureAD master application and setDataFromDataExcel
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AZUREAD")

# --- Grow the table (Tabla1) from 3 to 6 columns --------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F4")) | Out-Null

# --- New header cells (copy header style, then set text) -----------------
$ws.Range("A1:C1").Copy() | Out-Null
$ws.Range("D1:F1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Value = "URL"
$ws.Range("E1").Value = "CLIENT ID"
$ws.Range("F1").Value = "CLIENT SECRET"

# --- New 4th data row: Irisrusk master application ------------------------
$ws.Range("A4").Value = "AzureAD"
$ws.Range("B4").Value = "SGTO"
$ws.Range("C4").Value = "Irisrusk"
$ws.Range("D4").Value = "URL"
$ws.Range("E4").Value = "CLIENTID"
$ws.Range("F4").Value = "CLIENTSECRET"

# --- Placeholder "-" values for the pre-existing data rows ----------------
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "-"

$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "-"

# --- Widen the new columns -------------------------------------------------
$ws.Range("D1").EntireColumn.ColumnWidth = 29.83
$ws.Range("E1").EntireColumn.ColumnWidth = 29.83
$ws.Range("F1").EntireColumn.ColumnWidth = 29.83

# --- Make AZUREAD the active/selected sheet+cell --------------------------
$ws.Activate() | Out-Null
$ws.Range("F3").Select() | Out-Null
